# Revert capacity chart to show kilowatts (kW) instead of watts.
#  - Divide the "Energy Storage", "Solar" and "Wind" data-point values that were
#    entered in watts by 1000 so they read in kilowatts.
#  - Re-apply a one-decimal number format ("#,##0.0") to the data range so the
#    now-fractional kilowatt values still display sensibly.
#  - Rename the value-axis title from "Watts" to "Kilowatts (kW)" and simplify
#    its number format back to a plain "#,##0" (no more "...K" suffix logic,
#    since the axis is already in kW).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet data: convert watts -> kilowatts (divide by 1000) -----------

# Energy Storage (column C)
$ws.Range("C24").Value = 3.84
$ws.Range("C25").Value = 7.24
$ws.Range("C26").Value = 442

# Solar (column E)
$ws.Range("E12").Value = 18.4
$ws.Range("E13").Value = 38.7
$ws.Range("E14").Value = 21
$ws.Range("E15").Value = 32.4
$ws.Range("E17").Value = 10.7
$ws.Range("E18").Value = 23.6
$ws.Range("E19").Value = 50
$ws.Range("E20").Value = 32
$ws.Range("E21").Value = 146.8
$ws.Range("E22").Value = 123.7
$ws.Range("E23").Value = 266.945
$ws.Range("E24").Value = 338.97
$ws.Range("E25").Value = 329.85
$ws.Range("E26").Value = 357.71

# Wind (column G)
$ws.Range("G14").Value = 1.5

# The whole data block (B2:G26) keeps using the custom "#,##0" number format;
# now that kilowatt values can carry a fractional part, show one decimal.
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: axis title + axis number format ---------------------------------

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)  # xlValue
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
